$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample/test data (rows 2-4) with real book records.
# Written column-by-column so the shared-string table is built in the
# same bookname / accesscode / rfidcode column order as the source data.
$ws.Range("A2").Value = "毛概"
$ws.Range("A3").Value = "邓论"
$ws.Range("A4").Value = "旅游"

$ws.Range("B2").Value = "TP123"
$ws.Range("B3").Value = "TP12324"
$ws.Range("B4").Value = "A23132"

$ws.Range("C2").Value = "0x213213"
$ws.Range("C3").Value = "0x2343232"
$ws.Range("C4").Value = "0x1243245"

# Widen column C so the longer rfid hex codes are fully visible
$ws.Columns.Item(3).ColumnWidth = 16.43

# Move the active selection to C5
$ws.Range("C5").Select()
